# Append new rows of bitcoin news sentiment data (data updated until 09.09.2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data rows: Date (serial), column B, column C, column D
$data = @(
    @(45156, 4, 0, 0),
    @(45157, 1, 0, 1),
    @(45158, 2, 0, 0),
    @(45159, 4, 0, 0),
    @(45160, 2, 0, 1),
    @(45161, 3, 0, 0),
    @(45162, 3, 0, 0),
    @(45163, 1, 0, 0),
    @(45164, 1, 0, 0),
    @(45165, 2, 0, 0),
    @(45166, 2, 0, 0),
    @(45167, 5, 0, 0),
    @(45168, 4, 0, 1),
    @(45169, 4, 0, 1),
    @(45170, 3, 0, 0),
    @(45171, 3, 0, 0),
    @(45172, 3, 0, 0),
    @(45173, 2, 0, 1),
    @(45174, 2, 0, 0),
    @(45175, 9, 0, 0),
    @(45176, 4, 0, 1),
    @(45177, 2, 0, 0),
    @(45178, 3, 0, 0)
)

$startRow = 1292
$lastExistingRow = $startRow - 1

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    # Copy the formatting of the last existing data row down to the new row
    # so the new cells pick up the same style (date format in column A, etc.)
    $ws.Range("A" + $lastExistingRow + ":D" + $lastExistingRow).Copy() | Out-Null
    $ws.Range("A" + $r + ":D" + $r).PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
